$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.892.23"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.638.15"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.25"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5089"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2588"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06435"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.36"
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07795"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "1.661.95"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.274"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").Value = "1.865.64"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5600"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "0.0₅7677"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.31"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "25.901.04"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.88"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.390"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.952"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.156"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.789"
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.10"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1228"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.842"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.242"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04972"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.303"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.247"
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.387"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9047"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5565"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "1.137.53"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01575"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.64"
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.478"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8032"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "0.0₈113"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4249"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.806"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05067"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("E51").Value = "  +0.30%  "
